$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AddCustomerTest"

# Fill in header row
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# Fill in data row
$ws.Range("A2").Value = "Raman"
$ws.Range("B2").Value = "Arora"
$ws.Range("C2").Value = "A234wd"

# Set the selected cell as shown in the diff (selection C9)
$ws.Range("C9").Select()
